$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated AIP projection results after adding proactive auto-enrollment logic
# Row 2
$ws.Range("C2").Value = 9405
$ws.Range("D2").Value = 8333
$ws.Range("E2").Value = 0.8860180754917597
$ws.Range("F2").Value = 0.8840441332484616
$ws.Range("G2").Value = 0.1003492139685587
$ws.Range("H2").Value = 0.08871313388499895
$ws.Range("I2").Value = 42764307.097377
$ws.Range("J2").Value = 14897441.3715245
$ws.Range("L2").Value = 14897441.3715245
$ws.Range("M2").Value = 57661748.4689015
$ws.Range("N2").Value = 798735599.8913001
$ws.Range("O2").Value = 781035792.8873
$ws.Range("P2").Value = 0.01865128006508273
$ws.Range("Q2").Value = 0.01907395474982301

# Row 3
$ws.Range("C3").Value = 9584
$ws.Range("D3").Value = 8485
$ws.Range("E3").Value = 0.8853297161936561
$ws.Range("F3").Value = 0.882658899407053
$ws.Range("G3").Value = 0.1042274602239246
$ws.Range("H3").Value = 0.09199729532924165
$ws.Range("I3").Value = 51074341.32902661
$ws.Range("J3").Value = 18511671.10224956
$ws.Range("L3").Value = 18511671.10224956
$ws.Range("M3").Value = 69586012.43127617
$ws.Range("N3").Value = 845171543.037244
$ws.Range("O3").Value = 827576107.02422
$ws.Range("P3").Value = 0.02190285659136752
$ws.Range("Q3").Value = 0.02236854223451836

# Row 4
$ws.Range("C4").Value = 9768
$ws.Range("D4").Value = 8627
$ws.Range("E4").Value = 0.8831900081900081
$ws.Range("F4").Value = 0.8799469604243166
$ws.Range("G4").Value = 0.1088976469224528
$ws.Range("H4").Value = 0.09582415340677275
$ws.Range("I4").Value = 59723949.76822361
$ws.Range("J4").Value = 21601006.9797773
$ws.Range("L4").Value = 21601006.9797773
$ws.Range("M4").Value = 81324956.74800092
$ws.Range("N4").Value = 889821499.5132644
$ws.Range("O4").Value = 872257583.7186878
$ws.Range("P4").Value = 0.02427566314321819
$ws.Range("Q4").Value = 0.02476448171156726

# Row 5
$ws.Range("C5").Value = 9972
$ws.Range("D5").Value = 8834
$ws.Range("E5").Value = 0.8858804653028479
$ws.Range("F5").Value = 0.8834883488348835
$ws.Range("G5").Value = 0.1126647045506
$ws.Range("H5").Value = 0.09953795379537955
$ws.Range("I5").Value = 68799121.34358832
$ws.Range("J5").Value = 24398944.62936984
$ws.Range("L5").Value = 24398944.62936984
$ws.Range("M5").Value = 93198065.97295816
$ws.Range("N5").Value = 935647705.506358
$ws.Range("O5").Value = 918046483.2211035
$ws.Range("P5").Value = 0.02607706349919974
$ws.Range("Q5").Value = 0.02657702531985362

# Row 6
$ws.Range("C6").Value = 10163
$ws.Range("D6").Value = 9012
$ws.Range("E6").Value = 0.8867460395552494
$ws.Range("F6").Value = 0.8837026868013336
$ws.Range("G6").Value = 0.1144274300932091
$ws.Range("H6").Value = 0.1011198274171406
$ws.Range("I6").Value = 76969728.48235811
$ws.Range("J6").Value = 25719545.75328371
$ws.Range("L6").Value = 25719545.75328371
$ws.Range("M6").Value = 102689274.2356418
$ws.Range("N6").Value = 980150039.1428678
$ws.Range("O6").Value = 962442492.9427364
$ws.Range("P6").Value = 0.02624041700368162
$ws.Range("Q6").Value = 0.02672320262444395
